$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tracker")

# Row 10 - Catboost Test run
$ws.Range("A10").Value = 44971.973611111112
$ws.Range("B10").Value = "Test"
$ws.Range("D10").Value = "Catboost"
$ws.Range("E10").Value = "euclidean_dist, linear_dist, mean_hillshade, morning_hillshade, mean_amenties, aspect_dir, climatic_zone, geologic_zone, soil_type, scaling"
$ws.Range("F10").Value = "False, False, False, False, False, False, False, False, False"
$ws.Range("H10").Value = 0.95599999999999996
$ws.Range("I10").Value = 0.86499999999999999
$ws.Range("K10").Value = "Maria"

# Row 11 - Catboost Submission run
$ws.Range("A11").Value = 44971.974999999999
$ws.Range("B11").Value = "Submission"
$ws.Range("D11").Value = "Catboost"
$ws.Range("E11").Value = "euclidean_dist, linear_dist, mean_hillshade, morning_hillshade, mean_amenties, aspect_dir, climatic_zone, geologic_zone, soil_type, scaling"
$ws.Range("F11").Value = "False, False, False, False, False, False, False, False, False"
$ws.Range("H11").Value = 0.95
$ws.Range("J11").Value = 0.55500000000000005
$ws.Range("K11").Value = "Maria"

# Update the active cell selection to J12, matching the saved workbook state
$ws.Range("J12").Select()
